$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (regcntr_id, machine_id) pairs.
# Columns: A=regcntr_id B=machine_id C=lang_code D=is_active E=cr_by F=cr_dtimes G=eff_dtimes
$newRows = @(
    @(10002, 10021),
    @(10003, 10022),
    @(10004, 10023),
    @(10005, 10024),
    @(10006, 10025),
    @(10007, 10026),
    @(10008, 10027),
    @(10009, 10028),
    @(10010, 10029)
)

$startRow = 22
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $regcntrId = $newRows[$i][0]
    $machineId = $newRows[$i][1]

    $ws.Cells.Item($r, 1).Value = $regcntrId
    $ws.Cells.Item($r, 2).Value = $machineId
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Update selection / view to match the new active area.
$ws.Range("B22:B30").Select()

# Set page setup (portrait, matches new pageSetup element)
$ws.PageSetup.Orientation = 1
